# Update wilke validation results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("C2").Value = 111.4
$ws.Range("C3").Value = 305.60000000000002
$ws.Range("C4").Value = 472.2
$ws.Range("C5").Value = 582.70000000000005
$ws.Range("C6").Value = 1216.8
$ws.Range("C7").Value = 1247.2
$ws.Range("C8").Value = 2528.4
$ws.Range("C9").Value = 2159.1999999999998

$ws.Range("J20").Select()
